$d = $word.ActiveDocument

# --- Change 1: meeting length "30" -> "10" -------------------------------
# ("...attend this 30 ish minute meeting..." -> "...attend this 10 ish...")
$found1 = $d.Content.Find.Execute(
    "attend this 30 ", $true, $false, $false, $false, $false,
    $true, 1, $false, "attend this 10 ", 2)

# --- Change 2: first due date "June 18" -> "Oct 15th" --------------------
# ("(by Saturday June 18 11:59pm / Unit 8)" -> "(by Saturday Oct 15th  11:59pm / Unit 8)")
$found2 = $d.Content.Find.Execute(
    "Saturday June 18 11:59pm / Unit 8)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Saturday Oct 15th  11:59pm / Unit 8)", 2)

# Make the "th" in "15th" superscript, matching the vertAlign run the diff adds.
$rngTh = $d.Content
$foundTh = $rngTh.Find.Execute(
    "15th", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
[void]$rngTh.MoveStart(1, 2)
$rngTh.Font.Superscript = $true

# --- Change 3: second due date "June 25th " -> "Oct 22nd" ----------------
# ("(by Saturday June 25th  11:59pm / Unit 9)" -> "(by Saturday Oct 22nd 11:59pm / Unit 9)")
$found3 = $d.Content.Find.Execute(
    "June 25", $true, $false, $false, $false, $false,
    $true, 1, $false, "Oct 22", 2)

# The old superscript "th" run plus the following plain-space run collapse
# into a single superscript "nd" run (the diff drops the extra space run).
$rngNd = $d.Content
$foundNd = $rngNd.Find.Execute(
    "th  11:59pm / Unit 9)", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$trimLen = $rngNd.Text.Length
[void]$rngNd.MoveEnd(1, -($trimLen - 3))
$rngNd.Text = "nd"
